$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.466.55'
$ws.Range('E2').Value = '  -0.90%  '
$ws.Range('D3').Value = '3.251.75'
$ws.Range('E3').Value = '  +3.17%  '
$ws.Range('E4').Value = '  -0.04%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '595.96'
$c.ClearFormats()
$ws.Range('E5').Value = '  -0.90%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '140.62'
$c.ClearFormats()
$ws.Range('E6').Value = '  -1.10%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.244.50'
$ws.Range('E8').Value = '  +3.15%  '
$ws.Range('E9').Value = '  -1.63%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.148'
$c.ClearFormats()
$ws.Range('E10').Value = '  -0.96%  '
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('E12').Value = '  -0.37%  '
$ws.Range('E13').Value = '  -2.98%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '34.39'
$c.ClearFormats()
$ws.Range('E14').Value = '  -1.51%  '
$ws.Range('D15').Value = '3.781.57'
$ws.Range('E15').Value = '  +3.04%  '
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').Value = '3.249.90'
$ws.Range('E17').Value = '  +3.20%  '
$ws.Range('D18').Value = '63.445.46'
$ws.Range('E18').Value = '  -0.93%  '
$ws.Range('E19').Value = '  -0.83%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '476.20'
$c.ClearFormats()
$ws.Range('E20').Value = '  -1.96%  '
$ws.Range('E21').Value = '  -3.03%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '0.732'
$c.ClearFormats()
$ws.Range('E22').Value = '  +3.00%  '
$ws.Range('E23').Value = '  +2.78%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '84.08'
$c.ClearFormats()
$ws.Range('E24').Value = '  -4.58%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '13.17'
$c.ClearFormats()
$ws.Range('E25').Value = '  -0.48%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('E27').Value = '  -0.88%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '7.25'
$c.ClearFormats()
$ws.Range('E28').Value = '  +3.72%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '8.11'
$c.ClearFormats()
$ws.Range('E29').Value = '  -1.56%  '
$ws.Range('E30').Value = '  +2.75%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '27.50'
$c.ClearFormats()
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('E32').Value = '  -0.12%  '
$ws.Range('E33').Value = '  -4.13%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '2.53'
$c.ClearFormats()
$ws.Range('E34').Value = '  -4.40%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '1.08'
$c.ClearFormats()
$ws.Range('E35').Value = '  -1.26%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '5.94'
$c.ClearFormats()
$ws.Range('E36').Value = '  -1.90%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '52.82'
$c.ClearFormats()
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('D38').Value = '0.0₃0711'
$ws.Range('E38').Value = '  -3.91%  '
$ws.Range('E39').Value = '  -1.26%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '423.91'
$c.ClearFormats()
$ws.Range('E40').Value = '  -2.13%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '8.38'
$c.ClearFormats()
$ws.Range('E41').Value = '  +0.15%  '
$ws.Range('D42').Value = '2.973.86'
$ws.Range('E42').Value = '  +2.16%  '
$ws.Range('E43').Value = '  -5.83%  '
$ws.Range('E44').Value = '  -8.07%  '
$ws.Range('E45').Value = '  +1.90%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '2.18'
$c.ClearFormats()
$ws.Range('E46').Value = '  -0.64%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '25.89'
$c.ClearFormats()
$ws.Range('E48').Value = '  +0.22%  '
$ws.Range('E49').Value = '  -3.07%  '
$ws.Range('E50').Value = '  -0.64%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '121.80'
$c.ClearFormats()
$ws.Range('E51').Value = '  +0.59%  '
